# "Add files via upload" / "Changed U1 to SOT23-5L"
#
# Row 7 of the BOM sheet holds the U1 entry (A=Comment, B=Designator,
# C=Footprint, D=JLCPCB Part #). U1's footprint is being changed from
# SC-70-5 to SOT23-5L, which means the part itself changes too (from
# SN74LVC1G07DCKR / C7830 to TI SN74LVC1G125DBVR / C23654).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value = "SOT23-5L"
$ws.Range("D7").Value = "C23654"
$ws.Range("A7").Value = "TI SN74LVC1G125DBVR"

# Carry the Footprint column's cell formatting onto the updated part-number
# cell, matching how the edit was originally made.
$ws.Range("C7").Copy()
$ws.Range("D7").PasteSpecial(-4122)

# Leave the selection where the author left it after editing.
$ws.Range("A7").Select()
